$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = 45964
$ws.Range("B68").Value = "四方坪站"
$ws.Range("C68").Value = 9515.19
$ws.Range("D68").Value = 8366.44
$ws.Range("E68").Value = 3163.72
$ws.Range("F68").Value = 403

$ws.Range("A69").Value = 45964
$ws.Range("B69").Value = "高岭站"
$ws.Range("C69").Value = 3710.22
$ws.Range("D69").Value = 3169.26
$ws.Range("E69").Value = 1016.93
$ws.Range("F69").Value = 137

$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I69").Select()
